# Handback status report: refresh the handoff/handback timestamps for the
# 6a3367c8-... file's Xliff generation/handoff/handback events.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
# "Latest HO Xliff Generate Date" for 6a3367c8-082f-4270-9565-35c1ef070aca.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 21:06:20"

# --- zh-cn sheet ------------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the
# zh-cn xliff of 6a3367c8-082f-4270-9565-35c1ef070aca.md
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 21:06:15"
$wsZhCn.Range("K2").Value = "2016-08-22 21:06:33"

# --- de-de sheet ------------------------------------------------------------
# "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the
# de-de xliff of 6a3367c8-082f-4270-9565-35c1ef070aca.md
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 21:06:20"
$wsDeDe.Range("K2").Value = "2016-08-22 21:06:40"
